# Updates cryptos list price/volume figures per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.951.47'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.408.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.91'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.09'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.61'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.46%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.840.78'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.856.96'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.404.73'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.18'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.52'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '325.82'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.77'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.55%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.58'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.39'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.79'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0767'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.10'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.10'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.53%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.39'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.22%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.19'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '322.90'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.58%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '146.25'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.49%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0962'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0514'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.575'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0221'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.05'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.939'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.07%  '
